$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing header strings: U1 stays pointing at the same logical
# column but its label text changes; V1's label text changes; and a new
# column W1 is introduced reusing the old "param_P_to_charging_station1"
# label.
$ws.Range("U1").Value = "param_P_to_charging_station2"
$ws.Range("V1").Value = "param_E_pv3_solar"

# New column W, header styled like the other header cells (bold, centered,
# bordered) by copying U1's format, then overwrite with the right text.
$ws.Range("U1").Copy($ws.Range("W1"))
$ws.Range("W1").Value = "param_P_to_charging_station1"

# Data rows 2..50: refresh column U ("charging station" draw, now under the
# new "2" label), set column V to the flat 0.12 tariff used everywhere, and
# populate the brand-new column W with the values that used to live in U
# (now broken out as its own "charging station 1" series).
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0.12
$ws.Range("W2").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0.12
$ws.Range("W3").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0.12
$ws.Range("W4").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0.12
$ws.Range("W5").Value = 0
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = 0.12
$ws.Range("W6").Value = 0
$ws.Range("U7").Value = 15.84000000000001
$ws.Range("V7").Value = 0.12
$ws.Range("W7").Value = 26.60808333333334
$ws.Range("U8").Value = 72.63712500000004
$ws.Range("V8").Value = 0.12
$ws.Range("W8").Value = 67.93600000000002
$ws.Range("U9").Value = 101.3138194444444
$ws.Range("V9").Value = 0.12
$ws.Range("W9").Value = 75.17216666666668
$ws.Range("U10").Value = 0
$ws.Range("V10").Value = 0.12
$ws.Range("W10").Value = 0
$ws.Range("U11").Value = 0
$ws.Range("V11").Value = 0.12
$ws.Range("W11").Value = 0
$ws.Range("U12").Value = 0
$ws.Range("V12").Value = 0.12
$ws.Range("W12").Value = 144.3644583333334
$ws.Range("U13").Value = 106.8478888888889
$ws.Range("V13").Value = 0.12
$ws.Range("W13").Value = 33.96204166666667
$ws.Range("U14").Value = 22.44000000000001
$ws.Range("V14").Value = 0.12
$ws.Range("W14").Value = 0
$ws.Range("U15").Value = 38.50000000000002
$ws.Range("V15").Value = 0.12
$ws.Range("W15").Value = 0
$ws.Range("U16").Value = 75.26169444444447
$ws.Range("V16").Value = 0.12
$ws.Range("W16").Value = 0
$ws.Range("U17").Value = 67.76000000000001
$ws.Range("V17").Value = 0.12
$ws.Range("W17").Value = 136.6308472222223
$ws.Range("U18").Value = 31.84148611111112
$ws.Range("V18").Value = 0.12
$ws.Range("W18").Value = 77.6505277777778
$ws.Range("U19").Value = 0
$ws.Range("V19").Value = 0.12
$ws.Range("W19").Value = 95.36648611111116
$ws.Range("U20").Value = 0
$ws.Range("V20").Value = 0.12
$ws.Range("W20").Value = 0
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = 0.12
$ws.Range("W21").Value = 0
$ws.Range("U22").Value = 0
$ws.Range("V22").Value = 0.12
$ws.Range("W22").Value = 0
$ws.Range("U23").Value = 0
$ws.Range("V23").Value = 0.12
$ws.Range("W23").Value = 0
$ws.Range("U24").Value = 0
$ws.Range("V24").Value = 0.12
$ws.Range("W24").Value = 0
$ws.Range("U25").Value = 0
$ws.Range("V25").Value = 0.12
$ws.Range("W25").Value = 0
$ws.Range("U26").Value = 0
$ws.Range("V26").Value = 0.12
$ws.Range("W26").Value = 0
$ws.Range("U27").Value = 0
$ws.Range("V27").Value = 0.12
$ws.Range("W27").Value = 0
$ws.Range("U28").Value = 0
$ws.Range("V28").Value = 0.12
$ws.Range("W28").Value = 0
$ws.Range("U29").Value = 0
$ws.Range("V29").Value = 0.12
$ws.Range("W29").Value = 0
$ws.Range("U30").Value = 0
$ws.Range("V30").Value = 0.12
$ws.Range("W30").Value = 0
$ws.Range("U31").Value = 0
$ws.Range("V31").Value = 0.12
$ws.Range("W31").Value = 0
$ws.Range("U32").Value = 88.8306527777778
$ws.Range("V32").Value = 0.12
$ws.Range("W32").Value = 157.0561666666667
$ws.Range("U33").Value = 103.3040555555556
$ws.Range("V33").Value = 0.12
$ws.Range("W33").Value = 95.70000000000002
$ws.Range("U34").Value = 0
$ws.Range("V34").Value = 0.12
$ws.Range("W34").Value = 66.5776527777778
$ws.Range("U35").Value = 0
$ws.Range("V35").Value = 0.12
$ws.Range("W35").Value = 0
$ws.Range("U36").Value = 52.96148611111114
$ws.Range("V36").Value = 0.12
$ws.Range("W36").Value = 125.345
$ws.Range("U37").Value = 93.50488888888891
$ws.Range("V37").Value = 0.12
$ws.Range("W37").Value = 70.40000000000002
$ws.Range("U38").Value = 0
$ws.Range("V38").Value = 0.12
$ws.Range("W38").Value = 0
$ws.Range("U39").Value = 0
$ws.Range("V39").Value = 0.12
$ws.Range("W39").Value = 0
$ws.Range("U40").Value = 88.50600000000003
$ws.Range("V40").Value = 0.12
$ws.Range("W40").Value = 0
$ws.Range("U41").Value = 67.26270833333335
$ws.Range("V41").Value = 0.12
$ws.Range("W41").Value = 62.45250000000002
$ws.Range("U42").Value = 0
$ws.Range("V42").Value = 0.12
$ws.Range("W42").Value = 0
$ws.Range("U43").Value = 0
$ws.Range("V43").Value = 0.12
$ws.Range("W43").Value = 188.156375
$ws.Range("U44").Value = 0
$ws.Range("V44").Value = 0.12
$ws.Range("W44").Value = 14.52000000000001
$ws.Range("U45").Value = 0
$ws.Range("V45").Value = 0.12
$ws.Range("W45").Value = 0
$ws.Range("U46").Value = 0
$ws.Range("V46").Value = 0.12
$ws.Range("W46").Value = 0
$ws.Range("U47").Value = 0
$ws.Range("V47").Value = 0.12
$ws.Range("W47").Value = 0
$ws.Range("U48").Value = 0
$ws.Range("V48").Value = 0.12
$ws.Range("W48").Value = 0
$ws.Range("U49").Value = 0
$ws.Range("V49").Value = 0.12
$ws.Range("W49").Value = 0
$ws.Range("U50").Value = 0
$ws.Range("V50").Value = 0.12
$ws.Range("W50").Value = 0
